# edit.ps1 - apply the author's changes to Tool_Resources/drawings.pptx
#
# 1) Bump the cached "datetimeFigureOut" footer date from 7/11/16 to 7/12/16
#    everywhere it is cached (the slide master's Date Placeholder and each
#    of the 11 slide layouts' Date Placeholder).
# 2) On slide 6 ("SMOOT" poster textbox), shorten the long placeholder
#    sentence down to "... available" (keeping the leading "No Pareto
#    front " run intact), which also lets the autosize textbox shrink to
#    fit the now-shorter text.

$p = $ppt.ActivePresentation
$newDate = "7/12/16"

function Find-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            return $sh
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) Update the cached date field text ("7/11/16" -> "7/12/16").
# ---------------------------------------------------------------------

# 1a) Slide master's Date Placeholder.
$master = $p.Slides.Item(1).Master
(Find-DatePlaceholder $master.Shapes).TextFrame.TextRange.Text = $newDate

# 1b) Every slide layout's Date Placeholder.
#     `Master.CustomLayouts.Item(n)` only reliably resolves layout 1 in
#     this host, so instead we spin up one throwaway slide per layout
#     (via the classic ppLayout enum, which *does* resolve distinctly),
#     reach the layout through `Slide.CustomLayout`, edit it there, then
#     remove the scratch slides again (removing the slide does not
#     remove the shared layout part, so the edit sticks).
$layoutEnums = @(1, 2, 33, 4, 5, 6, 7, 8, 9, 10, 27)

$scratchSlides = @()
foreach ($layoutEnum in $layoutEnums) {
    $scratch = $p.Slides.Add($p.Slides.Count + 1, $layoutEnum)
    $layoutShapes = $scratch.CustomLayout.Shapes
    (Find-DatePlaceholder $layoutShapes).TextFrame.TextRange.Text = $newDate
    $scratchSlides += $scratch
}

# Remove the scratch slides again, highest index first so earlier
# entries in the list stay valid while deleting.
for ($i = $scratchSlides.Count - 1; $i -ge 0; $i--) {
    $scratchSlides[$i].Delete()
}

# ---------------------------------------------------------------------
# 2) Slide 6: shorten the "No Pareto front ..." sentence to "... available".
# ---------------------------------------------------------------------

$posterSlide = $p.Slides.Item(6)
$posterShape = $posterSlide.Shapes.Item(1)
$posterText = $posterShape.TextFrame.TextRange

# Second paragraph holds the sentence that needs trimming.
$sentence = $posterText.Paragraphs(2, 1)

$prefix = "No Pareto front "
$tail = $sentence.Characters($prefix.Length + 1, $sentence.Length - $prefix.Length)
$tail.Text = "available"
